# Applies the cryptos.xlsx price/volume/coin-swap update described by the commit diff.
# Values are written via a helper that forces Text storage (NumberFormat "@")
# so numeric-looking strings like "1.00" / "6.60" / "0.0000250" keep their exact
# textual representation instead of being coerced to numbers, then clears the
# temporary formatting so the cell style stays identical to the original (unstyled).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.ClearFormats()
}

Set-TextValue "D2" "69.765.22"
Set-TextValue "E2" "  +5.65%  "
Set-TextValue "D3" "3.430.64"
Set-TextValue "E3" "  +11.94%  "
Set-TextValue "E4" "  -0.09%  "
Set-TextValue "D5" "589.38"
Set-TextValue "E5" "  +2.75%  "
Set-TextValue "D6" "185.44"
Set-TextValue "E6" "  +9.41%  "
Set-TextValue "D7" "0.997"
Set-TextValue "E7" "  -0.33%  "
Set-TextValue "D8" "3.440.25"
Set-TextValue "E8" "  +12.46%  "
Set-TextValue "D9" "0.533"
Set-TextValue "E9" "  +4.97%  "
Set-TextValue "D10" "6.60"
Set-TextValue "E10" "  +4.08%  "
Set-TextValue "E11" "  +5.61%  "
Set-TextValue "E12" "  +4.38%  "
Set-TextValue "B13" "Avalanche"
Set-TextValue "C13" "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue "D13" "38.53"
Set-TextValue "E13" "  +8.26%  "
Set-TextValue "B14" "ShibaInu"
Set-TextValue "C14" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D14" "0.0000250"
Set-TextValue "E14" "  +5.02%  "
Set-TextValue "D15" "3.992.80"
Set-TextValue "E15" "  +11.64%  "
Set-TextValue "D16" "69.684.40"
Set-TextValue "E16" "  +5.61%  "
Set-TextValue "E17" "  +1.37%  "
Set-TextValue "D18" "3.411.08"
Set-TextValue "E18" "  +11.24%  "
Set-TextValue "D19" "7.41"
Set-TextValue "E19" "  +7.01%  "
Set-TextValue "D20" "17.09"
Set-TextValue "E20" "  +2.16%  "
Set-TextValue "D21" "505.40"
Set-TextValue "E21" "  +4.03%  "
Set-TextValue "D22" "8.40"
Set-TextValue "E22" "  +9.18%  "
Set-TextValue "D23" "0.727"
Set-TextValue "E23" "  +6.05%  "
Set-TextValue "D24" "86.54"
Set-TextValue "E24" "  +4.91%  "
Set-TextValue "E25" "  +5.84%  "
Set-TextValue "E26" "  +9.09%  "
Set-TextValue "D27" "10.86"
Set-TextValue "E27" "  +7.09%  "
Set-TextValue "D28" "1.00"
Set-TextValue "E28" "  -0.05%  "
Set-TextValue "E29" "  +5.50%  "
Set-TextValue "E30" "  +11.77%  "
Set-TextValue "E31" "  +4.25%  "
Set-TextValue "D32" "29.84"
Set-TextValue "E32" "  +8.31%  "
Set-TextValue "D33" "0.0000104"
Set-TextValue "E33" "  +14.88%  "
Set-TextValue "E34" "  +4.80%  "
Set-TextValue "D35" "1.00"
Set-TextValue "E35" "  +0.10%  "
Set-TextValue "D36" "6.09"
Set-TextValue "E36" "  +9.29%  "
Set-TextValue "E37" "  +6.37%  "
Set-TextValue "D38" "49.09"
Set-TextValue "E38" "  +5.57%  "
Set-TextValue "D39" "0.331"
Set-TextValue "E39" "  +10.52%  "
Set-TextValue "D40" "2.10"
Set-TextValue "E40" "  +7.34%  "
Set-TextValue "D41" "0.130"
Set-TextValue "E41" "  +6.43%  "
Set-TextValue "D42" "50.14"
Set-TextValue "E42" "  +2.24%  "
Set-TextValue "D43" "8.69"
Set-TextValue "E43" "  +5.29%  "
Set-TextValue "D44" "2.83"
Set-TextValue "E44" "  +13.00%  "
Set-TextValue "D45" "415.48"
Set-TextValue "E45" "  +14.92%  "
Set-TextValue "D46" "2.939.40"
Set-TextValue "E46" "  +5.55%  "
Set-TextValue "B47" "InjectiveProtocol"
Set-TextValue "C47" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D47" "27.93"
Set-TextValue "E47" "  +15.00%  "
Set-TextValue "B48" "VeChain"
Set-TextValue "C48" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D48" "0.0360"
Set-TextValue "E48" "  +4.87%  "
Set-TextValue "D49" "135.55"
Set-TextValue "E49" "  +0.88%  "
Set-TextValue "E51" "  +13.86%  "
